# chore: adapt column header formatting to respective input file names
#
# The sheet contains a diff table whose header row used the generic
# suffixes "_old" / "_new". Rename those headers to use the concrete
# format-version suffixes "_FV2404" / "_FV2410" instead, then turn the
# data range into a proper Excel Table ("Table1") and freeze the header
# row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Determine the extent of the data -------------------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. Rename the header cells in row 1 -----------------------------------
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($null -eq $text) { continue }

    if ($text -like "*_old") {
        $cell.Value2 = ($text -replace "_old$", "_FV2404")
    } elseif ($text -like "*_new") {
        $cell.Value2 = ($text -replace "_new$", "_FV2410")
    }
}

# --- 2. Wrap the data range in an Excel Table -------------------------------
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Renamed headers, created Table1 over" $dataRange.Address() "and froze the header row."
